# Updated main GSC export data:
# Drop the oldest day's row (2025-11-05) from the "Chart" sheet's rolling
# date-indexed export. Excel's row delete shifts every subsequent row up
# by one, which re-aligns each remaining row with the date it already had
# (2025-11-06 moves into row 2, 2025-11-07 into row 3, etc.) and shrinks
# the used range from A1:D91 down to A1:D90 -- exactly what the export
# refresh produces when the reporting window rolls forward one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
